$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5676286220550537
$ws.Range("B1").Value = 1.167033910751343
$ws.Range("C1").Value = 5.639540672302246
$ws.Range("D1").Value = 3.290453195571899
$ws.Range("E1").Value = 1.249758005142212
